$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2 through 144) holds a "Förändrad" (changed/updated) date
# serial number that was bumped by one day (46060 -> 46061) for every row.
for ($row = 2; $row -le 144; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
